$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coding")

# Row 2 - Total en-route ANS costs / SES States (EU27+2)
$ws.Range("B2").Value = 6904
$ws.Range("C2").Value = 6979
$ws.Range("D2").Value = 7064
$ws.Range("E2").Value = 6852
$ws.Range("F2").Value = 6524
$ws.Range("G2").Value = 6672

# Row 3 - Other 10 States in the Route Charges System
$ws.Range("B3").Value = 1436
$ws.Range("C3").Value = 1504
$ws.Range("D3").Value = 1538
$ws.Range("E3").Value = 1574
$ws.Range("F3").Value = 1423
$ws.Range("G3").Value = 1572

# Row 4 - Total en-route ANS costs
$ws.Range("B4").Value = 8340
$ws.Range("C4").Value = 8483
$ws.Range("D4").Value = 8602
$ws.Range("E4").Value = 8426
$ws.Range("F4").Value = 7947
$ws.Range("G4").Value = 8244

# Row 5 - SES States (EU27+2) service units
$ws.Range("B5").Value = 115
$ws.Range("C5").Value = 122
$ws.Range("D5").Value = 125
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = 67
$ws.Range("G5").Value = 108

# Row 6 - Other 10 States service units
$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 36
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 33

# Row 7 - Total en-route service units (million TSUs)
$ws.Range("B7").Value = 148
$ws.Range("C7").Value = 157
$ws.Range("D7").Value = 161
$ws.Range("E7").Value = 68
$ws.Range("F7").Value = 87
$ws.Range("G7").Value = 141

# Row 8 - SES States (EU27+2) ANS costs per TSU
$ws.Range("B8").Value = 60
$ws.Range("C8").Value = 57
$ws.Range("D8").Value = 56
$ws.Range("E8").Value = 131
$ws.Range("F8").Value = 98
$ws.Range("G8").Value = 62

# Row 9 - Other 10 States ANS costs per TSU
$ws.Range("B9").Value = 44
$ws.Range("C9").Value = 43
$ws.Range("D9").Value = 42
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 72
$ws.Range("G9").Value = 48

# Row 10 - En-route ANS costs per TSU (Total)
$ws.Range("B10").Value = 56
$ws.Range("C10").Value = 54
$ws.Range("D10").Value = 53
$ws.Range("E10").Value = 123
$ws.Range("F10").Value = 92
$ws.Range("G10").Value = 58

# Update the active selection shown in the sheet view
$ws.Range("K8").Select()
